$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (bold + thin border) for new rows 20-22, column A, 
# matching the existing styled column A cells (copy format from A19).
$ws.Range("A19").Copy()
$ws.Range("A20:A22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 'How many luxury hotels are in Happy Valley ski resort'
$ws.Range("D2").Value = 'Happy Valley ski resort'
$ws.Range("G2").Value = 'tourism=hotel, stars=*'
$ws.Range("J2").Value = $null
$ws.Range("K2").Value = $null
$ws.Range("V2").Value = 'data queries'
$ws.Range("W2").Value = 1
# --- Row 3 ---
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 'What area are within 50 km from family physician services in Saskatchewan'
$ws.Range("D3").Value = 'Saskatchewan'
$ws.Range("E3").Value = 'Canada'
$ws.Range("G3").Value = 'amenity=doctor'
$ws.Range("J3").Value = 'Buffer'
$ws.Range("K3").Value = 'Overlay analysis'
$ws.Range("V3").Value = 'data queries,buffer,overlay analysis'
$ws.Range("W3").Value = 10
# --- Row 4 ---
$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 15
$ws.Range("C4").Value = 'What areas are inside 1000 foot of schools in El Cajon'
$ws.Range("D4").Value = 'El Cajon'
$ws.Range("G4").Value = 'amenity=school, amenity=kindergarten'
# --- Row 5 ---
$ws.Range("A5").Value = 18
$ws.Range("B5").Value = 22
$ws.Range("C5").Value = 'What areas are outside 150 meters from hospitals in Houston'
$ws.Range("G5").Value = 'amenity=hospital'
# --- Row 6 ---
$ws.Range("A6").Value = 27
$ws.Range("B6").Value = 31
$ws.Range("C6").Value = 'What areas are outside 60 meters from water body in Houston'
$ws.Range("D6").Value = 'Houston'
$ws.Range("G6").Value = 'landuse=aquaculture, basin, salt_pond'
# --- Row 7 ---
$ws.Range("A7").Value = 31
$ws.Range("B7").Value = 36
$ws.Range("C7").Value = 'What areas are within 1000 meters of roads in Assam '
$ws.Range("D7").Value = 'Assam'
$ws.Range("E7").Value = 'India'
$ws.Range("G7").Value = 'highway=*'
# --- Row 8 ---
$ws.Range("B8").Value = 39
$ws.Range("C8").Value = 'What areas are within 1000 meters of the schools in Oleander'
$ws.Range("D8").Value = 'Oleander'
$ws.Range("G8").Value = 'amenity=school'
# --- Row 9 ---
$ws.Range("A9").Value = 50
$ws.Range("B9").Value = 56
$ws.Range("C9").Value = 'What areas do have altitude between 700 and 2000 meters in Spain'
$ws.Range("D9").Value = 'Spain'
$ws.Range("G9").Value = $null
$ws.Range("I9").Value = 'Topography'
$ws.Range("J9").Value = 'classification'
$ws.Range("K9").Value = 'Data queries'
$ws.Range("L9").Value = 'Data model conversion'
$ws.Range("V9").Value = 'topography,classification,data queries,data model conversion,overlay analysis'
$ws.Range("W9").Value = 17
# --- Row 10 ---
$ws.Range("A10").Value = 62
$ws.Range("B10").Value = 73
$ws.Range("C10").Value = 'What is the average rating of street pavement for each borough in New York City'
$ws.Range("D10").Value = 'New York City'
$ws.Range("G10").Value = $null
$ws.Range("I10").Value = 'Data editing'
$ws.Range("J10").Value = 'Overlay analysis'
$ws.Range("K10").Value = 'Data editing'
$ws.Range("L10").Value = 'Data queries'
$ws.Range("V10").Value = 'data editing,overlay analysis,data editing,data queries'
$ws.Range("W10").Value = 26
# --- Row 11 ---
$ws.Range("A11").Value = 65
$ws.Range("B11").Value = 78
$ws.Range("C11").Value = 'What is the cervix cancer mortality rate of white females for each city in the Western USA from 1970 to 1994'
$ws.Range("D11").Value = 'the Western USA '
$ws.Range("F11").Value = 'from 1970 to 1994'
$ws.Range("I11").Value = 'Data editing'
$ws.Range("J11").Value = 'Data queries'
$ws.Range("K11").Value = $null
$ws.Range("L11").Value = $null
$ws.Range("M11").Value = $null
$ws.Range("V11").Value = 'data editing,data queries'
$ws.Range("W11").Value = 9
# --- Row 12 ---
$ws.Range("A12").Value = 70
$ws.Range("B12").Value = 87
$ws.Range("C12").Value = 'What is the Euclidean distance to recreational sites in Utrecht'
$ws.Range("D12").Value = 'Utrecht'
$ws.Range("G12").Value = 'landuse=recreation_ground'
# --- Row 13 ---
$ws.Range("A13").Value = 72
$ws.Range("B13").Value = 89
$ws.Range("C13").Value = 'What is the Euclidean distance to subway stations in Amsterdam'
$ws.Range("D13").Value = 'Amsterdam'
$ws.Range("G13").Value = $null
# --- Row 14 ---
$ws.Range("A14").Value = 82
$ws.Range("B14").Value = 99
$ws.Range("C14").Value = 'What is the mean center of customers weighted by the number of transactions in Oleander city'
$ws.Range("D14").Value = 'Oleander city'
$ws.Range("E14").Value = $null
# --- Row 15 ---
$ws.Range("A15").Value = 85
$ws.Range("B15").Value = 102
$ws.Range("C15").Value = 'What is the mean center of the fire calls weighted by the priority in Fort Worth'
$ws.Range("D15").Value = 'Fort Worth'
$ws.Range("G15").Value = $null
$ws.Range("J15").Value = 'Overlay analysis'
$ws.Range("K15").Value = 'Geostatistics  '
$ws.Range("V15").Value = 'data queries,overlay analysis,geostatistics  '
$ws.Range("W15").Value = 35
# --- Row 16 ---
$ws.Range("A16").Value = 88
$ws.Range("B16").Value = 106
$ws.Range("C16").Value = 'What liquor stores are within 1000 foot of libraries in El Cajon'
$ws.Range("D16").Value = 'El Cajon'
$ws.Range("G16").Value = 'shop=alcohol, amenity=library'
$ws.Range("J16").Value = 'buffer'
$ws.Range("K16").Value = 'Overlay analysis'
$ws.Range("L16").Value = 'Data queries'
$ws.Range("V16").Value = 'data queries,buffer,overlay analysis,data queries'
$ws.Range("W16").Value = 21
# --- Row 17 ---
$ws.Range("A17").Value = 101
# --- Row 18 ---
$ws.Range("A18").Value = 102
$ws.Range("B18").Value = 122
$ws.Range("C18").Value = 'Where are the ski pistes in Happy Valley ski resort'
$ws.Range("D18").Value = 'Happy Valley ski resort'
$ws.Range("G18").Value = 'site=piste'
$ws.Range("I18").Value = 'Data queries'
$ws.Range("J18").Value = 'Geometry measurement'
$ws.Range("K18").Value = 'Data queries'
$ws.Range("V18").Value = 'data queries,geometry measurement,data queries'
$ws.Range("W18").Value = 8
# --- Row 19 ---
$ws.Range("A19").Value = 104
$ws.Range("B19").Value = 125
$ws.Range("C19").Value = 'Which houses are within 2 minutes driving time from fire stations  (from my current location) in Oleander'
$ws.Range("D19").Value = 'Oleander'
$ws.Range("G19").Value = 'amenity=fire_station'
$ws.Range("J19").Value = 'Network analysis'
$ws.Range("K19").Value = 'classification'
$ws.Range("M19").Value = 'Overlay analysis'
$ws.Range("V19").Value = 'data queries,network analysis,classification,data queries,overlay analysis'
$ws.Range("W19").Value = 0
# --- Row 20 ---
$ws.Range("A20").Value = 106
$ws.Range("B20").Value = 127
$ws.Range("C20").Value = 'Which houses have construction year between 1990 and 2000 in Utrecht'
$ws.Range("D20").Value = 'Utrecht'
$ws.Range("G20").Value = 'year_of_construction=*'
$ws.Range("H20").Value = 'done'
$ws.Range("I20").Value = 'Data queries'
$ws.Range("V20").Value = 'data queries'
$ws.Range("W20").Value = 1
$ws.Range("X20").Value = $false
# --- Row 21 ---
$ws.Range("A21").Value = 107
$ws.Range("B21").Value = 128
$ws.Range("C21").Value = 'Which land use contains meteorological stations in Netherlands'
$ws.Range("D21").Value = 'Netherlands'
$ws.Range("G21").Value = ' man_made=monitoring_station'
$ws.Range("H21").Value = 'done'
$ws.Range("I21").Value = 'Data queries'
$ws.Range("J21").Value = 'Overlay analysis'
$ws.Range("K21").Value = 'Data queries'
$ws.Range("V21").Value = 'data queries,overlay analysis,data queries'
$ws.Range("W21").Value = 2
$ws.Range("X21").Value = $false
# --- Row 22 ---
$ws.Range("A22").Value = 113
$ws.Range("B22").Value = 135
$ws.Range("C22").Value = 'Which wind farm proposals are nearest to the roads in Scotland'
$ws.Range("D22").Value = 'Scotland'
$ws.Range("G22").Value = 'highway=*'
$ws.Range("H22").Value = 'done'
$ws.Range("I22").Value = 'Data queries'
$ws.Range("J22").Value = 'network analysis'
$ws.Range("K22").Value = 'Data queries'
$ws.Range("V22").Value = 'data queries,network analysis,data queries'
$ws.Range("W22").Value = 4
$ws.Range("X22").Value = $false
